# Update the metadata sheet to match the target revision:
#  - Row 1: human-readable (capitalized) column headers
#  - Row 2: iaest-measure: / sdmx-dimension: identifiers per column
#  - Row 3: "medida" (measure) vs "dim" (dimension) classification per column
#  - Row 4: XSD datatype per column (xsd:string / xsd:int / xsd:date)
#  - Row 5: B5 "mapping-ano.xlsx" moved to H5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "Edad (grupos quinquenales)"
$ws.Range("B1").Value = "Personas"
$ws.Range("C1").Value = "Residencia comarca código"
$ws.Range("D1").Value = "Residencia CCAA nombre"
$ws.Range("E1").Value = "Residencia comarca nombre"
$ws.Range("F1").Value = "Residencia provincia nombre"
$ws.Range("G1").Value = "Extranjeros"
$ws.Range("H1").Value = "Año"
$ws.Range("I1").Value = "Sexo"

# Row 2 - measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:edad-grupos-quinquenales"
$ws.Range("B2").Value = "iaest-measure:personas"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("E2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("F2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("G2").Value = "iaest-measure:extranjeros"
$ws.Range("H2").Value = "sdmx-dimension:refPeriod"
$ws.Range("I2").Value = "iaest-measure:sexo"

# Row 3 - measure ("medida") vs dimension ("dim") classification
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "dim"
$ws.Range("I3").Value = "medida"

# Row 4 - XSD datatype
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "xsd:string"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:date"
$ws.Range("I4").Value = "xsd:string"

# Row 5 - the mapping reference moves from B5 to H5.
# Drop the whole row first so the old B5 cell entry does not linger,
# then write the new H5 cell and give it the same style as its column.
$ws.Rows("5").Delete()
$ws.Range("H5").Value = "mapping-ano.xlsx"
$ws.Range("H4").Copy()
$ws.Range("H5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
